# Added code to take dynamic screenshots
# - Switches the "Execute" flags for a couple of test rows to FALSE
# - Switches the "Mode" for a couple of rows from "remote" to "local"
# - Leaves the workbook focused/selected on the "Count" sheet at B2,
#   while remembering that "TestData" was last viewed at F2.

$wb = $excel.ActiveWorkbook

$wsCount    = $wb.Worksheets.Item("Count")
$wsTestData = $wb.Worksheets.Item("TestData")

# --- Count sheet: flip the Execute flag for row 2 off ---
$wsCount.Range("B2").Value = $false

# --- TestData sheet: flip Execute flag for row 2 off, and switch Mode
#     from "remote" to "local" for rows 2 and 5 ---
$wsTestData.Range("F2").Value = $false
$wsTestData.Range("G2").Value = "local"
$wsTestData.Range("G5").Value = "local"

# --- Update per-sheet selection state (mirrors Excel recording the last
#     selected cell on each sheet, and which sheet tab is active) ---
$wsTestData.Activate()
$wsTestData.Range("F2").Select()

$wsCount.Activate()
$wsCount.Range("B2").Select()
